$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AW2").Value = 151.105868
$ws.Range("AK3").Value = 41.944178
$ws.Range("AN4").Value = 2.108993
$ws.Range("AQ5").Value = 13.89522
$ws.Range("AW6").Value = 106.051528
$ws.Range("AW7").Value = 108.838079
$ws.Range("AK8").Value = 8.943194
$ws.Range("AW9").Value = 143.875579
$ws.Range("AW10").Value = 92.91621499999999
$ws.Range("AW11").Value = 76.99950200000001
$ws.Range("AW12").Value = 77.901875
$ws.Range("AN13").Value = 20.798044
$ws.Range("AW14").Value = 155.959572
$ws.Range("AW15").Value = 141.774722
$ws.Range("AW16").Value = 110.903252
$ws.Range("AW17").Value = 113.200556
$ws.Range("AQ18").Value = 3.637025
$ws.Range("AN19").Value = 16.073588
$ws.Range("AW20").Value = 80.787384
$ws.Range("AW21").Value = 72.87614600000001
$ws.Range("AK22").Value = 8.943206
$ws.Range("AK23").Value = 8.943251999999999
$ws.Range("AN24").Value = 3.88728
$ws.Range("AW25").Value = 106.050995
$ws.Range("AN26").Value = 0.940544
$ws.Range("AK27").Value = 8.943356
$ws.Range("AW28").Value = 141.780231
$ws.Range("AW29").Value = 62.863715
$ws.Range("AW30").Value = 71.957199
$ws.Range("AW31").Value = 2.097245
$ws.Range("AW32").Value = 143.870486
$ws.Range("AW33").Value = 80.09920099999999
$ws.Range("AN34").Value = 0.940486
$ws.Range("AW35").Value = 0.878646
$ws.Range("AW36").Value = 108.145
$ws.Range("AW37").Value = 115.930741
$ws.Range("AW38").Value = 3.162639
$ws.Range("AK39").Value = 8.943229000000001
$ws.Range("AK40").Value = 8.943125
$ws.Range("AK41").Value = 8.943438
$ws.Range("AK42").Value = 8.943148000000001
$ws.Range("AQ43").Value = 13.895093
$ws.Range("AW44").Value = 65.182002
$ws.Range("AW45").Value = 80.09858800000001
$ws.Range("AW46").Value = 31.215359
$ws.Range("AW47").Value = 31.214144
$ws.Range("AK48").Value = 30.992743
$ws.Range("AW49").Value = 108.145984
$ws.Range("AW50").Value = 108.839491
$ws.Range("AW51").Value = 15.982106
$ws.Range("AW52").Value = 122.054977
$ws.Range("AW53").Value = 99.812512
$ws.Range("AW54").Value = 146.917164
$ws.Range("AW55").Value = 136.965984
$ws.Range("AK56").Value = 31.792986
$ws.Range("AW57").Value = 108.838657
$ws.Range("AW58").Value = 3.166748
$ws.Range("AQ59").Value = 13.894606
$ws.Range("AW60").Value = 44.840231
$ws.Range("AW61").Value = 112.09434
$ws.Range("AW62").Value = 72.89548600000001
$ws.Range("AK63").Value = 3.968958
$ws.Range("AN64").Value = 20.79706
$ws.Range("AW65").Value = 148.210301
$ws.Range("AW66").Value = 151.103194
$ws.Range("AW67").Value = 108.854433
$ws.Range("AW68").Value = 65.955961
$ws.Range("AW69").Value = 8.988519
$ws.Range("AK70").Value = 8.943414000000001
$ws.Range("AK71").Value = 8.943171
$ws.Range("AW72").Value = 1.992558
$ws.Range("AK73").Value = 30.992824
$ws.Range("AW74").Value = 122.925868
$ws.Range("AW75").Value = 59.952083
$ws.Range("AN76").Value = 20.802488
$ws.Range("AW77").Value = 64.160729
$ws.Range("AN78").Value = 2.108727
$ws.Range("AW79").Value = 136.969375
$ws.Range("AW80").Value = 31.215139
$ws.Range("AW81").Value = 16.130498
$ws.Range("AW82").Value = 10.825255
$ws.Range("AW83").Value = 0.8773030000000001
$ws.Range("AW84").Value = 115.994699
$ws.Range("AW85").Value = 146.869444
$ws.Range("AW86").Value = 94.83910899999999
$ws.Range("AW87").Value = 64.15901599999999
$ws.Range("AW88").Value = 44.833079
$ws.Range("AQ89").Value = 28.924884
$ws.Range("AW90").Value = 100.84162
$ws.Range("AW91").Value = 141.77434
$ws.Range("AW92").Value = 51.214769
$ws.Range("AK93").Value = 3.966516
$ws.Range("AW94").Value = 108.837882
$ws.Range("AK95").Value = 8.943078999999999
$ws.Range("AW96").Value = 148.210694
$ws.Range("AW97").Value = 59.825359
$ws.Range("AQ98").Value = 13.895324
$ws.Range("AW99").Value = 77.894537
$ws.Range("AW100").Value = 69.970394
$ws.Range("AW101").Value = 143.960301
$ws.Range("AW102").Value = 99.957292
$ws.Range("AW103").Value = 39.216042
$ws.Range("AK104").Value = 41.944144
$ws.Range("AK105").Value = 3.968981
$ws.Range("AK106").Value = 31.793044
$ws.Range("AW107").Value = 98.394167
$ws.Range("AW108").Value = 107.037685
$ws.Range("AW109").Value = 59.881308
$ws.Range("AW110").Value = 8.970000000000001
$ws.Range("AK111").Value = 8.943113
$ws.Range("AK112").Value = 8.943299
$ws.Range("AN113").Value = 2.108681
$ws.Range("AW114").Value = 151.106157
$ws.Range("AW115").Value = 136.917037
$ws.Range("AW116").Value = 146.867581
$ws.Range("AW117").Value = 39.214444
$ws.Range("AW118").Value = 10.815891
$ws.Range("AN119").Value = 0.940532
$ws.Range("AW120").Value = 125.83625
$ws.Range("AW121").Value = 99.97730300000001
$ws.Range("AW122").Value = 77.84774299999999
$ws.Range("AW123").Value = 64.156863
$ws.Range("AK124").Value = 8.943241
$ws.Range("AK125").Value = 31.79309
$ws.Range("AW126").Value = 141.774954
$ws.Range("AK127").Value = 8.943148000000001
$ws.Range("AN128").Value = 2.10875
$ws.Range("AW129").Value = 136.869549
$ws.Range("AK130").Value = 30.992859
$ws.Range("AK131").Value = 8.943125
$ws.Range("AW132").Value = 51.077118
$ws.Range("AW133").Value = 45.22853
$ws.Range("AW134").Value = 154.919572
$ws.Range("AW135").Value = 80.75993099999999
$ws.Range("AW136").Value = 39.2139
$ws.Range("AK137").Value = 8.943218
$ws.Range("AN138").Value = 15.998461
$ws.Range("AW139").Value = 136.920243
$ws.Range("AW140").Value = 141.776609
$ws.Range("AW141").Value = 77.957234
$ws.Range("AK142").Value = 8.943160000000001
$ws.Range("AW143").Value = 112.093912
$ws.Range("AW144").Value = 108.83838
$ws.Range("AW145").Value = 71.958229
$ws.Range("AK146").Value = 8.943403
$ws.Range("AW147").Value = 44.090405
$ws.Range("AW148").Value = 43.890463
$ws.Range("AQ149").Value = 28.924769
$ws.Range("AW150").Value = 65.7886
$ws.Range("AW151").Value = 55.853831
$ws.Range("AW152").Value = 85.92400499999999
$ws.Range("AW153").Value = 39.213206
$ws.Range("AK154").Value = 8.943263999999999
$ws.Range("AQ155").Value = 20.904884
$ws.Range("AW156").Value = 108.145787
$ws.Range("AK157").Value = 59.995185
$ws.Range("AW158").Value = 64.162477
$ws.Range("AK159").Value = 3.966516
$ws.Range("AW160").Value = 98.1686
$ws.Range("AW161").Value = 13.915046
$ws.Range("AK162").Value = 8.943403
$ws.Range("AW163").Value = 146.870382
$ws.Range("AW164").Value = 123.220683
$ws.Range("AW165").Value = 62.863507
$ws.Range("AK166").Value = 8.94309
$ws.Range("AW167").Value = 98.173542
$ws.Range("AW168").Value = 151.106551
$ws.Range("AW169").Value = 72.874618
$ws.Range("AW170").Value = 143.869931
$ws.Range("AW171").Value = 134.976771
$ws.Range("AW172").Value = 72.89537
$ws.Range("AW173").Value = 1.992396
$ws.Range("AK174").Value = 30.992731
$ws.Range("AW175").Value = 143.869572
$ws.Range("AW176").Value = 99.973113
$ws.Range("AW177").Value = 45.227813
$ws.Range("AW178").Value = 77.78752299999999
$ws.Range("AN179").Value = 0.940463
$ws.Range("AW180").Value = 151.11412
$ws.Range("AW181").Value = 85.93619200000001
$ws.Range("AW182").Value = 64.159988
$ws.Range("AW183").Value = 45.231354
$ws.Range("AW184").Value = 99.976933
$ws.Range("AK185").Value = 8.943241
$ws.Range("AK186").Value = 30.992859
$ws.Range("AW187").Value = 143.960579
$ws.Range("AW188").Value = 158.181412
$ws.Range("AW189").Value = 64.16583300000001
$ws.Range("AW190").Value = 37.953229
$ws.Range("AW191").Value = 131.197118
$ws.Range("AW192").Value = 141.775521
$ws.Range("AW193").Value = 69.879074
$ws.Range("AW194").Value = 125.83309
$ws.Range("AW195").Value = 73.80251199999999
$ws.Range("AK196").Value = 3.96897
$ws.Range("AW197").Value = 146.843715
$ws.Range("AW198").Value = 143.960035
$ws.Range("AW199").Value = 142.882257
$ws.Range("AW200").Value = 154.864838
$ws.Range("AW201").Value = 157.005231
$ws.Range("AW202").Value = 146.918414
$ws.Range("AW203").Value = 108.840104
$ws.Range("AW204").Value = 14.183113
$ws.Range("AW205").Value = 108.858079
$ws.Range("AK206").Value = 30.992743
$ws.Range("AW207").Value = 155.960556
$ws.Range("AW208").Value = 146.859479
$ws.Range("AW209").Value = 155.956424
$ws.Range("AW210").Value = 56.910729
$ws.Range("AK211").Value = 8.943194
$ws.Range("AK212").Value = 8.943182999999999
$ws.Range("AQ213").Value = 13.894826
$ws.Range("AN214").Value = 0.9404400000000001
$ws.Range("AW215").Value = 64.16518499999999
$ws.Range("AW216").Value = 115.93066
$ws.Range("AW217").Value = 133.171968
$ws.Range("AW218").Value = 72.941574
$ws.Range("AW219").Value = 65.181771
$ws.Range("AW220").Value = 45.229664
$ws.Range("AK221").Value = 1.884965
$ws.Range("AW222").Value = 108.146319
$ws.Range("AW223").Value = 64.02478000000001
$ws.Range("AK224").Value = 8.943333000000001
$ws.Range("AK225").Value = 30.992847
$ws.Range("AQ226").Value = 20.897627
$ws.Range("AW227").Value = 98.172315
$ws.Range("AW228").Value = 143.844815
$ws.Range("AW229").Value = 108.838889
$ws.Range("AW230").Value = 72.872963
$ws.Range("AW231").Value = 16.127546
$ws.Range("AW232").Value = 44.932431
$ws.Range("AW233").Value = 39.212326
$ws.Range("AK234").Value = 8.943137
$ws.Range("AK235").Value = 31.793102
$ws.Range("AW236").Value = 108.839294
$ws.Range("AK237").Value = 8.943218
$ws.Range("AW238").Value = 80.09982599999999
$ws.Range("AW239").Value = 128.992488
$ws.Range("AK240").Value = 8.936215000000001
$ws.Range("AK241").Value = 8.94309
$ws.Range("AW242").Value = 108.106898
$ws.Range("AW243").Value = 39.222072
$ws.Range("AW244").Value = 39.219919
$ws.Range("AW245").Value = 13.804468
$ws.Range("AK246").Value = 8.943171
$ws.Range("AW247").Value = 106.051157
$ws.Range("AW248").Value = 30.210637
$ws.Range("AW249").Value = 15.982859
$ws.Range("AW250").Value = 136.969595
$ws.Range("AW251").Value = 146.917674
$ws.Range("AW252").Value = 125.833438
$ws.Range("AW253").Value = 15.983669
$ws.Range("AN254").Value = 2.108704
$ws.Range("AW255").Value = 1.992292
$ws.Range("AW256").Value = 125.833981
$ws.Range("AW257").Value = 45.224063
$ws.Range("AT258").Value = 1.131157
$ws.Range("AW259").Value = 135.862789
$ws.Range("AW260").Value = 51.07831
$ws.Range("AK261").Value = 8.943368
$ws.Range("AW262").Value = 3.164919
$ws.Range("AW263").Value = 122.05478
$ws.Range("AW264").Value = 2.097975
$ws.Range("AQ265").Value = 13.894734
$ws.Range("AW266").Value = 31.214063
$ws.Range("AK267").Value = 8.943379999999999
$ws.Range("AW268").Value = 123.221053
$ws.Range("AW269").Value = 134.97772
$ws.Range("AW270").Value = 98.17337999999999
$ws.Range("AW271").Value = 77.901447
$ws.Range("AW272").Value = 19.228866
$ws.Range("AW273").Value = 122.183449
$ws.Range("AK274").Value = 30.992755
$ws.Range("AT275").Value = 1.130938
$ws.Range("AW276").Value = 136.18265
$ws.Range("AW277").Value = 80.790544
$ws.Range("AW278").Value = 141.773553
$ws.Range("AW279").Value = 64.1664
$ws.Range("AW280").Value = 108.841227
$ws.Range("AK281").Value = 60.009039
$ws.Range("AW282").Value = 133.172431
$ws.Range("AW283").Value = 64.159271
$ws.Range("AW284").Value = 38.970741
$ws.Range("AK285").Value = 8.943345000000001
$ws.Range("AW286").Value = 73.932211
$ws.Range("AW287").Value = 148.212153
$ws.Range("AW288").Value = 37.884363
$ws.Range("AK289").Value = 8.943414000000001
$ws.Range("AK290").Value = 8.943102
$ws.Range("AN291").Value = 0.940509
$ws.Range("AW292").Value = 37.882928
$ws.Range("AW293").Value = 37.884699
$ws.Range("AK294").Value = 31.793102
$ws.Range("AW295").Value = 141.773935
$ws.Range("AW296").Value = 94.83562499999999
$ws.Range("AW297").Value = 13.915995
$ws.Range("AW298").Value = 94.838414
$ws.Range("AW299").Value = 108.143657
$ws.Range("AK300").Value = 3.966505
$ws.Range("AW301").Value = 151.107002
$ws.Range("AW302").Value = 45.223206
$ws.Range("AW303").Value = 8.933287
$ws.Range("AN304").Value = 22.813125
$ws.Range("AW305").Value = 129.881053
$ws.Range("AK306").Value = 8.943322
$ws.Range("AK307").Value = 8.943275
$ws.Range("AW308").Value = 84.94006899999999
$ws.Range("AW309").Value = 108.840556
$ws.Range("AW310").Value = 77.847442
$ws.Range("AW311").Value = 23.85816
$ws.Range("AW312").Value = 3.161134
$ws.Range("AN313").Value = 20.79625
$ws.Range("AW314").Value = 87.86503500000001
$ws.Range("AK315").Value = 8.943345000000001
$ws.Range("AW316").Value = 13.915509
$ws.Range("AW317").Value = 80.09647
$ws.Range("AW318").Value = 155.965382
$ws.Range("AW319").Value = 66.16435199999999
$ws.Range("AW320").Value = 45.230567
$ws.Range("AN321").Value = 1.90647
$ws.Range("AW322").Value = 113.203727
$ws.Range("AK323").Value = 8.943322
$ws.Range("AK324").Value = 8.943275
$ws.Range("AK325").Value = 8.943368
$ws.Range("AW326").Value = 69.924109
$ws.Range("AW327").Value = 48.862627
$ws.Range("AW328").Value = 13.91662
$ws.Range("AN329").Value = 20.797292
$ws.Range("AW330").Value = 122.183669
$ws.Range("AW331").Value = 51.077859
$ws.Range("AW332").Value = 39.211563
$ws.Range("AW333").Value = 45.284074
$ws.Range("AK334").Value = 8.943391
$ws.Range("AK335").Value = 31.793044
$ws.Range("AW336").Value = 146.868287
$ws.Range("AW337").Value = 148.210463
$ws.Range("AK338").Value = 8.943287
$ws.Range("AK339").Value = 16.894722
$ws.Range("AW340").Value = 154.859132
$ws.Range("AW341").Value = 45.231782
$ws.Range("AW342").Value = 125.834502
$ws.Range("AW343").Value = 16.129051
$ws.Range("AW344").Value = 14.164086
$ws.Range("AW345").Value = 141.773287
$ws.Range("AW346").Value = 113.201042
$ws.Range("AW347").Value = 31.801238
$ws.Range("AW348").Value = 115.9936
$ws.Range("AW349").Value = 64.160625
$ws.Range("AW350").Value = 150.999039
$ws.Range("AW351").Value = 146.904896
$ws.Range("AN352").Value = 15.691157
$ws.Range("AW353").Value = 148.211262
$ws.Range("AW354").Value = 125.835023
$ws.Range("AW355").Value = 155.956782
$ws.Range("AW356").Value = 3.89316
$ws.Range("AW357").Value = 141.771817
$ws.Range("AW358").Value = 108.837697
$ws.Range("AW359").Value = 99.97032400000001
$ws.Range("AW360").Value = 16.128368
$ws.Range("AK361").Value = 8.94331
$ws.Range("AK362").Value = 30.992836
$ws.Range("AW363").Value = 62.863854
$ws.Range("AW364").Value = 62.863391
$ws.Range("AK365").Value = 8.943426000000001
$ws.Range("AQ366").Value = 13.894977
$ws.Range("AN367").Value = 1.906539
